$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  "zögern","wärmen","helfen","nennen","kümmern","lügen","zählen","messen","boxen","töten",
  "rasen","stehlen","zeigen","lenken","treiben","graben","weichen","fischen","scheinen","schnellen",
  "schwingen","spüren","gelten","klingeln","werden","hacken","rechnen","brauchen","trotzen","bauen",
  "zünden","kranken","flehen","kratzen","sinken","tollen","dringen","suchen","nutzen","machen",
  "achten","decken","schulden","wetten","planen","bremsen","quälen","bitten","räumen","binden",
  "lockern","runden","küssen","sprechen","kehren","lösen","arten","loben","leiden","wachsen",
  "zielen","fassen","enden","reiten","schauen","drucken","fließen","stechen","spinnen","flüchten",
  "schrecken","regnen","tropfen","sichern","streichen","heißen","malen","liefern","sperren","schenken",
  "gründen","schätzen","leisten","schmecken","wundern","fahren","beißen","mühen","hören","stimmen",
  "fangen","reizen","drehen","münzen","klingen","sitzen","spielen","kosten","saufen","ärgern",
  "sorgen","parken","schicken","lesen","fallen","trennen","löschen","ehren","wollen","werfen",
  "rufen","kichern","dienen","schreiten","zahlen","grüßen","pfeifen","fällen","fesseln","wehtun",
  "ändern","albern","steuern","schlagen","mögen","wirken","bergen","weigern","feiern","schlucken",
  "siegen","bluten","heben","herzen","formen","spannen","seufzen","hauen","wüten","biegen",
  "greifen","wenden","ahnen","äußern","warnen","bellen","stammen","liegen","schleppen","heulen",
  "pflegen","betteln","reisen","schwören","schwächen","rücken","morden","sprengen","filmen","tauschen",
  "trauen","sterben","irren","stecken","freuen","folgen","schließen","bieten","narren","erben",
  "weinen","altern","trösten","mauern","starren","führen","öffnen","klettern","scheitern","backen",
  "geben","dürfen","jagen","jubeln","treffen","füllen","leugnen","knarren","pflanzen","platzen",
  "posten","heilen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $words[$i]
}
